$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old table content ---------------------------------
$ws.Range("A1:C5").ClearContents()

# Remove the two existing hyperlinks (surveys.healthvibe.eu / google.be)
$ws.Range("A1").Hyperlinks.Delete()

# --- Re-enter the new data ---------------------------------------------
# NOTE: the order in which new text values are typed in determines the
# order they end up in the workbook's shared string table, so these
# assignments are intentionally ordered to reproduce that table exactly.
$ws.Range("B2").Value2 = "google.be"
$ws.Range("C2").Value2 = "https://google.be"
$ws.Range("B1").Value2 = "title"
$ws.Range("B3").Value2 = "facebook"
$ws.Range("C3").Value2 = "https://facebook.be"
$ws.Range("A2").Value2 = "REQ"
$ws.Range("A1").Value2 = "VAL"
$ws.Range("C1").Value2 = "siteURL"
$ws.Range("D1").Value2 = "justavalue"
$ws.Range("A3").Value2 = "REQ"

# --- Hyperlinks for the URL column --------------------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "https://google.be")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://facebook.be")

# Re-apply the Hyperlink cell style (the Hyperlinks.Add call above tends
# to create a brand new style record instead of reusing the workbook's
# existing built-in "Hyperlink" style, so set it explicitly afterwards).
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"

# C5 keeps the hyperlink-style formatting even though it no longer holds
# a value or a hyperlink itself.
$ws.Range("C5").Style = "Hyperlink"

# --- Misc view state -----------------------------------------------------
$ws.Range("E13").Select() | Out-Null
